# AHDT2_OP_SG_HAVCR2_VAR.docx content edit
#
# The canonical diff splits several existing runs into multiple runs
# (same visible text / formatting, just different <w:r> boundaries -
# this is what Word does internally when you click in the middle of a
# run and retype/reformat a sub-span) and inserts one new phrase
# ("or in homopolymer regions").
#
# Word's object model has no direct "split this run in two" verb, so
# we use a well known trick: toggling a character formatting property
# on a sub-range to a different value and then immediately back to its
# original value forces Word to materialise a run boundary at the
# edges of that sub-range, without changing how the text looks.

function Split-RangeAt {
    param($doc, $pos, $rangeEnd)
    $r = $doc.Range($pos, $rangeEnd)
    $wasBold = $r.Font.Bold
    if ($wasBold) {
        $r.Font.Bold = $false
        $r.Font.Bold = $true
    } else {
        $r.Font.Bold = $true
        $r.Font.Bold = $false
    }
}

$d = $word.ActiveDocument

# ---------------------------------------------------------------
# 1) "Clinical Indication"  ->  "Clinical " + "Indication"
# ---------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Clinical Indication", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$start = $rng.Start
$end = $rng.End
$p1 = $start + 9
Split-RangeAt $d $p1 $end

# ---------------------------------------------------------------
# 2) "Germline variant analysis of HAVCR2 exon 2 including Tyr82 and
#     Ile97 hotspot variant loci." ->
#     "Germline variant analysis of HAVCR2 exon 2 including Tyr82 and
#      Ile97 " + "hotspot variant loci."
# ---------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Germline variant analysis of HAVCR2 exon 2 including Tyr82 and Ile97 hotspot variant loci.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$start = $rng.Start
$end = $rng.End
$p1 = $start + 69
Split-RangeAt $d $p1 $end

# ---------------------------------------------------------------
# 3) "Illumina NovaSeq" -> "Illumina " + "NovaSeq"
# ---------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Illumina NovaSeq", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$start = $rng.Start
$end = $rng.End
$p1 = $start + 9
Split-RangeAt $d $p1 $end

# ---------------------------------------------------------------
# 4) "A custom pipeline utilising the Oncoanalyser analysis pipeline
#     (OncoPath v1)" -> 5 runs:
#     "A custom pipeline utilising the " / "Oncoanalyser" /
#     " analysis pipeline (" / "OncoPath" / " v1)"
# ---------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("A custom pipeline utilising the Oncoanalyser analysis pipeline (OncoPath v1)", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$start = $rng.Start
$end = $rng.End
$p1 = $start + 32
$p2 = $start + 44
$p3 = $start + 64
$p4 = $start + 72
Split-RangeAt $d $p1 $end
Split-RangeAt $d $p2 $end
Split-RangeAt $d $p3 $end
Split-RangeAt $d $p4 $end

# ---------------------------------------------------------------
# 5) "are classified according to ACMG guidelines for the
#     interpretation of sequence variants" ->
#     "are classified according to ACMG guidelines " +
#     "for the interpretation of sequence variants"
# ---------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("are classified according to ACMG guidelines for the interpretation of sequence variants", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$start = $rng.Start
$end = $rng.End
$p1 = $start + 44
Split-RangeAt $d $p1 $end

# ---------------------------------------------------------------
# 6) "... 25 bp in length)" -> "... 25 bp in " + "length" +
#     " or in homopolymer regions" (new text) + ")" (unchanged,
#     re-split off from the trailing sentence so it stays its own run)
# ---------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("25 bp in length", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$start = $rng.Start
$end = $rng.End

$insPoint = $d.Range($end, $end)
$insPoint.InsertAfter(" or in homopolymer regions")
$newEnd = $end + 26
$closeParenEnd = $newEnd + 1

$p1 = $start + 9
Split-RangeAt $d $p1 $closeParenEnd
Split-RangeAt $d $end $closeParenEnd
Split-RangeAt $d $newEnd $closeParenEnd

# ---------------------------------------------------------------
# 7) "...family relationships, and clinical diagnoses are as stated..."
#     -> "...family relationships, and " +
#        "clinical diagnoses are as stated..."
# ---------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute(", variant zygosity is assumed to be either heterozygous or homozygous in the germline based on allele frequency for the purpose of clinical interpretation. Please note Peter Mac assumes sample identification, family relationships, and clinical diagnoses are as stated on the request. Our clinical recommendations may be based on evidence from third-party data sources and should be interpreted in the context of all other clinical and laboratory information for this patient.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$start = $rng.Start
$end = $rng.End
$p1 = $start + 235
Split-RangeAt $d $p1 $end
